$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.395797515623563
$ws.Range("D2").Value = 7.916171316735013
$ws.Range("E2").Value = 12.82347712581739
$ws.Range("F2").Value = 37.78070605950429
$ws.Range("G2").Value = 43.27792541234198
$ws.Range("H2").Value = 17.79846701943315
$ws.Range("I2").Value = 21.54351685093415
$ws.Range("J2").Value = 10.03009224842337
$ws.Range("K2").Value = 13.96414524196522
$ws.Range("M2").Value = 17.16354082990321
$ws.Range("N2").Value = 19.87187874832656

$ws.Range("B3").Value = 5.290501482306739
$ws.Range("D3").Value = 7.907760066427351
$ws.Range("E3").Value = 12.83401007021344
$ws.Range("F3").Value = 37.78104809818578
$ws.Range("G3").Value = 43.2070750696508
$ws.Range("H3").Value = 17.8356155312612
$ws.Range("I3").Value = 21.63384017979175
$ws.Range("J3").Value = 10.05179593171212
$ws.Range("K3").Value = 13.69884088830049
$ws.Range("M3").Value = 17.07279580918332
$ws.Range("N3").Value = 19.93655472237162

$ws.Range("B4").Value = 5.22570963216414
$ws.Range("D4").Value = 7.903549758240919
$ws.Range("E4").Value = 12.84237841870403
$ws.Range("F4").Value = 37.790566151704
$ws.Range("G4").Value = 43.17660695822381
$ws.Range("H4").Value = 17.86176092696136
$ws.Range("I4").Value = 21.69250092665276
$ws.Range("J4").Value = 10.06621231766221
$ws.Range("K4").Value = 13.53640856350783
$ws.Range("M4").Value = 17.02010126913667
$ws.Range("N4").Value = 19.97807633803815

$ws.Range("B5").Value = 5.199308668363292
$ws.Range("D5").Value = 7.902075457226503
$ws.Range("E5").Value = 12.84626703216499
$ws.Range("F5").Value = 37.79678409909848
$ws.Range("G5").Value = 43.16747323298054
$ws.Range("H5").Value = 17.87325289167611
$ws.Range("I5").Value = 21.7172118460362
$ws.Range("J5").Value = 10.07236156144383
$ws.Range("K5").Value = 13.4704260133897
$ws.Range("M5").Value = 16.99940490697095
$ws.Range("N5").Value = 19.99545339471487

$ws.Range("B6").Value = 5.194926051821601
$ws.Range("D6").Value = 7.901845273945303
$ws.Range("E6").Value = 12.84694164000367
$ws.Range("F6").Value = 37.79795780421536
$ws.Range("G6").Value = 43.16615489637577
$ws.Range("H6").Value = 17.87521166431017
$ws.Range("I6").Value = 21.72136379868042
$ws.Range("J6").Value = 10.07339922358669
$ws.Range("K6").Value = 13.45948507796738
$ws.Range("M6").Value = 16.99601569774345
$ws.Range("N6").Value = 19.99836646541931

$ws.Range("B7").Value = 5.225353520471842
$ws.Range("D7").Value = 7.903528895855934
$ws.Range("E7").Value = 12.84242892424647
$ws.Range("F7").Value = 37.79064054038624
$ws.Range("G7").Value = 43.1764704842333
$ws.Range("H7").Value = 17.86191252253711
$ws.Range("I7").Value = 21.6928309211218
$ws.Range("J7").Value = 10.06629413688191
$ws.Range("K7").Value = 13.53551772616615
$ws.Range("M7").Value = 17.01981898289102
$ws.Range("N7").Value = 19.97830884019794

$ws.Range("B8").Value = 5.359542565692582
$ws.Range("D8").Value = 7.913073986182567
$ws.Range("E8").Value = 12.82671458545462
$ws.Range("F8").Value = 37.77889172502474
$ws.Range("G8").Value = 43.25079479131792
$ws.Range("H8").Value = 17.81058263291044
$ws.Range("I8").Value = 21.57399642340144
$ws.Range("J8").Value = 10.03734958086465
$ws.Range("K8").Value = 13.87262693989371
$ws.Range("M8").Value = 17.13163490523961
$ws.Range("N8").Value = 19.89380416940566

$ws.Range("B9").Value = 5.619987556212911
$ws.Range("D9").Value = 7.939294618764425
$ws.Range("E9").Value = 12.81096251469837
$ws.Range("F9").Value = 37.82970938157337
$ws.Range("G9").Value = 43.49960583896998
$ws.Range("H9").Value = 17.73645294475263
$ws.Range("I9").Value = 21.36631919632111
$ws.Range("J9").Value = 9.989226621386663
$ws.Range("K9").Value = 14.53339623531386
$ws.Range("M9").Value = 17.37410541556003
$ws.Range("N9").Value = 19.74238758501389

$ws.Range("B10").Value = 5.807648926691398
$ws.Range("D10").Value = 7.963035624952786
$ws.Range("E10").Value = 12.80854087829081
$ws.Range("F10").Value = 37.91198963245863
$ws.Range("G10").Value = 43.74453142273373
$ws.Range("H10").Value = 17.69823589791529
$ws.Range("I10").Value = 21.22912884849794
$ws.Range("J10").Value = 9.959118107562185
$ws.Range("K10").Value = 15.01350207853104
$ws.Range("M10").Value = 17.56527227192494
$ws.Range("N10").Value = 19.63976182052728

$ws.Range("B11").Value = 5.891818367488371
$ws.Range("D11").Value = 7.974785246172286
$ws.Range("E11").Value = 12.80941736573726
$ws.Range("F11").Value = 37.95913226507827
$ws.Range("G11").Value = 43.8692214554391
$ws.Range("H11").Value = 17.68439084124353
$ws.Range("I11").Value = 21.17004525393207
$ws.Range("J11").Value = 9.946556580284485
$ws.Range("K11").Value = 15.22971847344746
$ws.Range("M11").Value = 17.65481106642963
$ws.Range("N11").Value = 19.59492597019015

$ws.Range("B12").Value = 5.923487653642112
$ws.Range("D12").Value = 7.97936895570904
$ws.Range("E12").Value = 12.81003272596484
$ws.Range("F12").Value = 37.97837393013333
$ws.Range("G12").Value = 43.91832257888874
$ws.Range("H12").Value = 17.67965795347401
$ws.Range("I12").Value = 21.14814903376196
$ws.Range("J12").Value = 9.94196277188683
$ws.Range("K12").Value = 15.31119893585976
$ws.Range("M12").Value = 17.68906485058046
$ws.Range("N12").Value = 19.57821218278128

$ws.Range("B13").Value = 5.916676653918078
$ws.Range("D13").Value = 7.978375831043184
$ws.Range("E13").Value = 12.80988760761771
$ws.Range("F13").Value = 37.97416820963778
$ws.Range("G13").Value = 43.90766442478684
$ws.Range("H13").Value = 17.68065457374274
$ws.Range("I13").Value = 21.15284355194398
$ws.Range("J13").Value = 9.942944886444028
$ws.Range("K13").Value = 15.29366952470171
$ws.Range("M13").Value = 17.68167260809834
$ws.Range("N13").Value = 19.58180004986155

$ws.Range("B14").Value = 5.894428077891669
$ws.Range("D14").Value = 7.975159670369248
$ws.Range("E14").Value = 12.80946231710077
$ws.Range("F14").Value = 37.96068748055092
$ws.Range("G14").Value = 43.87322343581175
$ws.Range("H14").Value = 17.68399123940044
$ws.Range("I14").Value = 21.16823427354134
$ws.Range("J14").Value = 9.946175380132138
$ws.Range("K14").Value = 15.23643030819007
$ws.Range("M14").Value = 17.65762231774092
$ws.Range("N14").Value = 19.5935456219038

$ws.Range("B15").Value = 5.880772732570686
$ws.Range("D15").Value = 7.973207113072226
$ws.Range("E15").Value = 12.80923869675411
$ws.Range("F15").Value = 37.95261090194971
$ws.Range("G15").Value = 43.85237185808627
$ws.Range("H15").Value = 17.68610147537609
$ws.Range("I15").Value = 21.17772369427112
$ws.Range("J15").Value = 9.948175367076962
$ws.Range("K15").Value = 15.20131571833896
$ws.Range("M15").Value = 17.64293536670354
$ws.Range("N15").Value = 19.60077454035132

$ws.Range("B16").Value = 5.802121400963552
$ws.Range("D16").Value = 7.962286685936967
$ws.Range("E16").Value = 12.8085233199758
$ws.Range("F16").Value = 37.90910364224008
$ws.Range("G16").Value = 43.73664766111004
$ws.Range("H16").Value = 17.69921201869165
$ws.Range("I16").Value = 21.23305694906855
$ws.Range("J16").Value = 9.95996185120765
$ws.Range("K16").Value = 14.99932075940225
$ws.Range("M16").Value = 17.55947058096217
$ws.Range("N16").Value = 19.6427290443478

$ws.Range("B17").Value = 5.753541351680461
$ws.Range("D17").Value = 7.955829082018131
$ws.Range("E17").Value = 12.80859042972407
$ws.Range("F17").Value = 37.88489662215419
$ws.Range("G17").Value = 43.66903734439474
$ws.Range("H17").Value = 17.70816228227103
$ws.Range("I17").Value = 21.26785322538029
$ws.Range("J17").Value = 9.96748300386532
$ws.Range("K17").Value = 14.87478400267532
$ws.Range("M17").Value = 17.50891121601882
$ws.Range("N17").Value = 19.66893943242027

$ws.Range("B18").Value = 5.725488334951478
$ws.Range("D18").Value = 7.952204368605746
$ws.Range("E18").Value = 12.80881528511361
$ws.Range("F18").Value = 37.87188812543618
$ws.Range("G18").Value = 43.63140061640745
$ws.Range("H18").Value = 17.71364340238705
$ws.Range("I18").Value = 21.28818012417985
$ws.Range("J18").Value = 9.971915816641665
$ws.Range("K18").Value = 14.80295259455449
$ws.Range("M18").Value = 17.48007469465001
$ws.Range("N18").Value = 19.68418908133886

$ws.Range("B19").Value = 5.715971981567895
$ws.Range("D19").Value = 7.950992544559919
$ws.Range("E19").Value = 12.80892343517982
$ws.Range("F19").Value = 37.8676409655497
$ws.Range("G19").Value = 43.61887301805631
$ws.Range("H19").Value = 17.71555640771224
$ws.Range("I19").Value = 21.29511624223161
$ws.Range("J19").Value = 9.97343504888263
$ws.Range("K19").Value = 14.77859972364258
$ws.Range("M19").Value = 17.47035372447721
$ws.Range("N19").Value = 19.68938229748868

$ws.Range("B20").Value = 5.758724508962015
$ws.Range("D20").Value = 7.956507254509831
$ws.Range("E20").Value = 12.80856401584185
$ws.Range("F20").Value = 37.88737887375317
$ws.Range("G20").Value = 43.6761052796961
$ws.Range("H20").Value = 17.70717502320164
$ws.Range("I20").Value = 21.26411671291686
$ws.Range("J20").Value = 9.966671308037975
$ws.Range("K20").Value = 14.88806257754539
$ws.Range("M20").Value = 17.5142682702918
$ws.Range("N20").Value = 19.66613127999472

$ws.Range("B21").Value = 5.900968798086256
$ws.Range("D21").Value = 7.976100705541217
$ws.Range("E21").Value = 12.80957955093036
$ws.Range("F21").Value = 37.96460944236153
$ws.Range("G21").Value = 43.88328867974023
$ws.Range("H21").Value = 17.68299733411534
$ws.Range("I21").Value = 21.16370069272779
$ws.Range("J21").Value = 9.945222084803204
$ws.Range("K21").Value = 15.25325422923051
$ws.Range("M21").Value = 17.66467722950546
$ws.Range("N21").Value = 19.5900884926774

$ws.Range("B22").Value = 5.992732185901944
$ws.Range("D22").Value = 7.989688358177006
$ws.Range("E22").Value = 12.81189489691798
$ws.Range("F22").Value = 38.02318043763175
$ws.Range("G22").Value = 44.02966092552308
$ws.Range("H22").Value = 17.67016840219457
$ws.Range("I22").Value = 21.10085570478439
$ws.Range("J22").Value = 9.932153575845639
$ws.Range("K22").Value = 15.48958524599794
$ws.Range("M22").Value = 17.76499194475965
$ws.Range("N22").Value = 19.54193167859019

$ws.Range("B23").Value = 5.943876185100192
$ws.Range("D23").Value = 7.982365549831102
$ws.Range("E23").Value = 12.81050838902117
$ws.Range("F23").Value = 37.99118186903863
$ws.Range("G23").Value = 43.95054496363923
$ws.Range("H23").Value = 17.67674320365941
$ws.Range("I23").Value = 21.13414285449999
$ws.Range("J23").Value = 9.939041657917787
$ws.Range("K23").Value = 15.3636911136203
$ws.Range("M23").Value = 17.71127558762215
$ws.Range("N23").Value = 19.56749326084399

$ws.Range("B24").Value = 5.75638158644933
$ws.Range("D24").Value = 7.956200378937618
$ws.Range("E24").Value = 12.80857537727074
$ws.Range("F24").Value = 37.8862538174716
$ws.Range("G24").Value = 43.67290602240961
$ws.Range("H24").Value = 17.70762031799338
$ws.Range("I24").Value = 21.26580498693799
$ws.Range("J24").Value = 9.967037936790515
$ws.Range("K24").Value = 14.88206005506166
$ws.Range("M24").Value = 17.51184562557791
$ws.Range("N24").Value = 19.66740028192932

$ws.Range("B25").Value = 5.550035802011029
$ws.Range("D25").Value = 7.931407615515017
$ws.Range("E25").Value = 12.81361402385054
$ws.Range("F25").Value = 37.80805731783553
$ws.Range("G25").Value = 43.42132132872691
$ws.Range("H25").Value = 17.75365881997586
$ws.Range("I25").Value = 21.41979393563365
$ws.Range("J25").Value = 10.00132235010411
$ws.Range("K25").Value = 14.35521927587067
$ws.Range("M25").Value = 17.30613904213309
$ws.Range("N25").Value = 19.78182914108059
